$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 6, column F: "CRESS" -> "Unclassified"
$ws.Range("F6").Value = "Unclassified"

# Add a new data row (row 14) for the batCV-Sc703 / JN857329 record.
# Copy formatting from row 7 (an existing "Unclassified" styled row) first,
# then fill in the values.
$ws.Range("A7:H7").Copy()
$ws.Range("A14:H14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A14").Value = "JN857329"
$ws.Range("C14").Value = "Circoviridae batCV-SC703"
$ws.Range("E14").Value = "CRESS-4"
$ws.Range("B14").Value = "batCV-Sc703"
$ws.Range("D14").Value = "Circoviridae"
$ws.Range("F14").Value = "Unclassified"
$ws.Range("G14").Value = "Unknown"
$ws.Range("H14").Value = "Unknown"

# Update the active cell selection to reflect where the user ended up (C20).
$null = $ws.Range("C20").Select()
